$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the changed data values
$ws.Range("C2").Value = 12.2
$ws.Range("C3").Value = 11.1
$ws.Range("B4").Value = 0.7
$ws.Range("C4").Value = 1.65
$ws.Range("C5").Value = 28

# The last row (row 6) was a blank leftover row - remove it so the
# used range shrinks back to A1:C5
$ws.Rows.Item(6).Delete()

# Resize columns A and C (the ColumnWidth property is expressed in
# "characters" and gets snapped to the nearest whole pixel internally,
# so the inputs below are chosen to land as closely as possible on the
# target stored widths of 27 and 27.25 characters)
$ws.Columns.Item(1).ColumnWidth = 26.3
$ws.Columns.Item(3).ColumnWidth = 26.6

# Update the selected cell
[void]$ws.Range("C2").Select()

# Page setup: portrait, paper size 9 (A4)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
